$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "38.169.64"
$ws.Range("E2").Value = "  +3.11%  "
$ws.Range("D3").Value = "2.062.28"
$ws.Range("E3").Value = "  +2.69%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "230.91"
$ws.Range("E5").Value = "  +2.66%  "
$ws.Range("E6").Value = "  +2.91%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "58.61"
$ws.Range("E7").Value = "  +6.99%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("E9").Value = "  +3.13%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0811"
$ws.Range("E10").Value = "  +3.82%  "
$ws.Range("E11").Value = "  +0.83%  "
$ws.Range("D12").Value = "2.366.90"
$ws.Range("E12").Value = "  +2.77%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "14.69"
$ws.Range("E13").Value = "  +4.27%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "20.74"
$ws.Range("E14").Value = "  +2.91%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.755"
$ws.Range("E15").Value = "  +2.21%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.28"
$ws.Range("E16").Value = "  +3.53%  "
$ws.Range("D17").Value = "2.064.75"
$ws.Range("E17").Value = "  +2.17%  "
$ws.Range("D18").Value = "38.056.11"
$ws.Range("E18").Value = "  +3.01%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.15"
$ws.Range("E19").Value = "  +0.14%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "70.02"
$ws.Range("E20").Value = "  +2.12%  "
$ws.Range("D21").Value = "0.0₃0833"
$ws.Range("E21").Value = "  +2.42%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "225.15"
$ws.Range("E22").Value = "  +1.49%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.00"
$ws.Range("E23").Value = "  -0.04%  "
$ws.Range("E24").Value = "  +0.90%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.25"
$ws.Range("E25").Value = "  +3.60%  "
$ws.Range("E26").Value = "  +2.34%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "166.52"
$ws.Range("E27").Value = "  +0.50%  "
$ws.Range("E28").Value = "  +7.30%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.07"
$ws.Range("E29").Value = "  +2.25%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.35"
$ws.Range("E30").Value = "  +0.24%  "
$ws.Range("E31").Value = "  +2.36%  "
$ws.Range("E32").Value = "  +1.55%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.63"
$ws.Range("E33").Value = "  +5.37%  "
$ws.Range("E34").Value = "  +0.94%  "
$ws.Range("E35").Value = "  +7.96%  "
$ws.Range("E36").Value = "  +0.75%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.04"
$ws.Range("E37").Value = "  +15.26%  "
$ws.Range("E38").Value = "  +5.98%  "
$ws.Range("E39").Value = "  +0.25%  "
$ws.Range("E40").Value = "  +2.40%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "98.51"
$ws.Range("E41").Value = "  +4.09%  "
$ws.Range("D42").Value = "1.480.31"
$ws.Range("E42").Value = "  +0.30%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0946"
$ws.Range("E43").Value = "  +3.35%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "16.85"
$ws.Range("E44").Value = "  +3.79%  "
$ws.Range("E45").Value = "  +3.63%  "
$ws.Range("E46").Value = "  +1.00%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.06"
$ws.Range("E47").Value = "  +15.98%  "
$ws.Range("E48").Value = "  +1.56%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.97"
$ws.Range("E49").Value = "  +1.93%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.11"
$ws.Range("E50").Value = "  -0.52%  "
$ws.Range("D51").Value = "2.255.09"
$ws.Range("E51").Value = "  +2.99%  "
